# Risk evaluation workbook: three risks got their measured probability
# ("Risk Probability of Occurrence") revised after the final time-tracking
# pass, which shuffles their ranking in the descending Risk-Factor sort
# that's already applied to the A7:G12 table.
#   - "Product cannot be completed within the semester": 20% -> 3%
#   - "Underlying (web) technology changes...":           10% -> 9%
#   - "Product quality is buggy...":                       18% -> 15%

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value  = 0.03   # Product cannot be completed within the semester
$ws.Range("C9").Value  = 0.09   # Underlying (web) technology changes and disables product
$ws.Range("C10").Value = 0.15   # Product quality is buggy and user satisfaction is low

# Re-sort the risk table (A7:G12) descending by Risk Factor (column E),
# same condition already recorded on the sheet's sortState.
$tableRange = $ws.Range("A7:G12")
$sortKey = $ws.Range("E7:E12")
$tableRange.Sort($sortKey, 2, $null, $null, 1, $null, $null, 1)

# Re-affirm the selection and scroll the window up one row so the header
# is back in view (topLeftCell A8 -> A7).
$ws.Range("C10").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
